# CCGX-Modbus-TCP-register-list.xlsx
# Add unit-id mappings for the Cerbo GX ports on the "Unit ID mapping" sheet,
# and log the change on the "Document versions" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "Unit ID mapping" sheet: insert 4 new rows right after the existing
#    "EasySolar-II/Multiplus-II GX VE.Bus port (ttyS3)" row (row 7) and fill
#    them with the new Cerbo GX port mappings.
# ---------------------------------------------------------------------------
$wsUnit = $wb.Worksheets.Item("Unit ID mapping")

$wsUnit.Range("A8:A11").EntireRow.Insert()

$wsUnit.Range("A8").Value2 = 227
$wsUnit.Range("B8").Value2 = 276
$wsUnit.Range("C8").Value2 = "Cerbo GX VE.Bus port (ttyS4)"
$wsUnit.Rows.Item(8).RowHeight = 13.8

$wsUnit.Range("A9").Value2 = 226
$wsUnit.Range("B9").Value2 = 279
$wsUnit.Range("C9").Value2 = "Cerbo GX VE.Direct port 1 (ttyS7)"
$wsUnit.Rows.Item(9).RowHeight = 13.8

$wsUnit.Range("A10").Value2 = 224
$wsUnit.Range("B10").Value2 = 278
$wsUnit.Range("C10").Value2 = "Cerbo GX VE.Direct port 2 (ttyS6)"
$wsUnit.Rows.Item(10).RowHeight = 13.8

$wsUnit.Range("A11").Value2 = 223
$wsUnit.Range("B11").Value2 = 277
$wsUnit.Range("C11").Value2 = "Cerbo GX VE.Direct port 3 (ttyS5)"
$wsUnit.Rows.Item(11).RowHeight = 13.8

# Column C needed to grow a lot to fit the longer remark text.
$wsUnit.Columns.Item(3).ColumnWidth = 64.15

# Restore the view roughly to where it was (scrolled near the top of the
# newly inserted block, with C7 selected).
$wsUnit.Activate()
$wsUnit.Range("C7").Select()

# ---------------------------------------------------------------------------
# 2. "Document versions" sheet: append the Rev 27 changelog entry.
# ---------------------------------------------------------------------------
$wsDoc = $wb.Worksheets.Item("Document versions")

$wsDoc.Range("A68").Value2 = "Rev 27"
$wsDoc.Range("B68").Value2 = "Add mappings for Cerbo GX ports"

$wsDoc.Activate()
$wsDoc.Range("B69").Select()
